# Update Wnt5a-Fzd3 LR-pairs sheet with newly recomputed TPM statistics.
# The "Resolving-Mac" target cluster is dropped entirely (both rows that
# referenced it as the Target cluster), and all remaining sending/target
# cluster combinations get refreshed downstream metric values (columns
# E through T) reflecting the new TPM-derived computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows whose Target cluster is "Resolving-Mac"
# (row 11 = MuSCs -> Resolving-Mac, row 6 = FAPs -> Resolving-Mac).
# Deleting from the bottom up keeps earlier row numbers stable.
$ws.Rows("11").Delete()
$ws.Rows("6").Delete()

# After the deletions, the remaining 8 rows (now rows 2-9) keep their
# original relative order:
#   2: FAPs  -> ECs
#   3: FAPs  -> FAPs
#   4: FAPs  -> MuSCs
#   5: FAPs  -> Neutrophils
#   6: MuSCs -> ECs
#   7: MuSCs -> FAPs
#   8: MuSCs -> MuSCs
#   9: MuSCs -> Neutrophils
# Only columns E:T (the numeric statistics) need to be refreshed with the
# recomputed values; columns A:D (cluster/symbol labels) are unchanged.

$newData = @{
    2 = @(3, 1, 7.514794999999999, 22.544385, 0.977669497583861, 0.977669497583861, 3, 1, 0.2310223333333334, 0.6930670000000001, 0.1371162007804382, 0.1371162007804382, 1.736085475421667, 15.624769278795, 0.1340543271276189, 0.1340543271276189)
    3 = @(3, 1, 7.514794999999999, 22.544385, 0.977669497583861, 0.977669497583861, 3, 1, 0.4828523333333334, 1.448557, 0.2865821521640898, 0.2865821521640898, 3.628536300271667, 32.656826702445, 0.2801826287227673, 0.2801826287227673)
    4 = @(3, 1, 7.514794999999999, 22.544385, 0.977669497583861, 0.977669497583861, 3, 1, 0.9450603333333335, 2.835181, 0.5609114952015948, 0.5609114952015947, 7.101934667631667, 63.91741200868501, 0.5483860597027554, 0.5483860597027553)
    5 = @(3, 1, 7.514794999999999, 22.544385, 0.977669497583861, 0.977669497583861, 1, 0.3333333333333333, 0.02593033333333333, 0.077791, 0.01539015185387714, 0.01539015185387714, 0.1948611392816667, 1.753750253535, 0.01504648203071939, 0.01504648203071939)
    6 = @(2, 0.6666666666666666, 0.171642, 0.514926, 0.02233050241613897, 0.02233050241613898, 3, 1, 0.2310223333333334, 0.6930670000000001, 0.1371162007804382, 0.1371162007804382, 0.039653135338, 0.356878218042, 0.003061873652819372, 0.003061873652819373)
    7 = @(2, 0.6666666666666666, 0.171642, 0.514926, 0.02233050241613897, 0.02233050241613898, 3, 1, 0.4828523333333334, 1.448557, 0.2865821521640898, 0.2865821521640898, 0.082877740198, 0.7458996617820001, 0.006399523441322515, 0.006399523441322516)
    8 = @(2, 0.6666666666666666, 0.171642, 0.514926, 0.02233050241613897, 0.02233050241613898, 3, 1, 0.9450603333333335, 2.835181, 0.5609114952015948, 0.5609114952015947, 0.162212045734, 1.459908411606, 0.01252543549883934, 0.01252543549883934)
    9 = @(2, 0.6666666666666666, 0.171642, 0.514926, 0.02233050241613897, 0.02233050241613898, 1, 0.3333333333333333, 0.02593033333333333, 0.077791, 0.01539015185387714, 0.01539015185387714, 0.004450734273999999, 0.040056608466, 0.0003436698231577492, 0.0003436698231577493)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i   # column E is index 5
        $ws.Cells.Item([int]$r, $col).Value = $vals[$i]
    }
}
